$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feedback")

# 1. Insert a new column before column A (everything shifts right by one)
$ws.Columns.Item(1).Insert()

# 2. Insert a new row before row 3 (pushes nothing yet - sheet currently only has
#    2 data rows after the column insert, so this creates row 3 fresh)
$ws.Rows.Item(3).Insert()

# 3. Row 1 (headers)
$ws.Cells.Item(1,1).Value  = "Page Name"
$ws.Cells.Item(1,2).Value  = "Contact name"
$ws.Cells.Item(1,3).Value  = "Primary role"
$ws.Cells.Item(1,4).Value  = "Recognition level"
$ws.Cells.Item(1,5).Value  = "contact attributes"
$ws.Cells.Item(1,6).Value  = "institution"
$ws.Cells.Item(1,7).Value  = "institution type"
$ws.Cells.Item(1,8).Value  = "payer"
$ws.Cells.Item(1,9).Value  = "IDS"
$ws.Cells.Item(1,10).Value = "IDS Input"
$ws.Cells.Item(1,11).Value = "Date"
$ws.Cells.Item(1,12).Value = "time spent"
$ws.Cells.Item(1,13).Value = "What topics were discussed?"

# 4. Row 2
$ws.Cells.Item(2,1).Value  = "Feedback"
$ws.Cells.Item(2,2).Value  = "Auto testing 6"
$ws.Cells.Item(2,3).Value  = "Medical Director"
$ws.Cells.Item(2,4).Value  = "National"
$ws.Cells.Item(2,5).Value  = "Cardiology"
$ws.Cells.Item(2,6).Value  = "Auto Testing 4"
$ws.Cells.Item(2,7).Value  = "Academic"
$ws.Cells.Item(2,8).Value  = "yes"
$ws.Cells.Item(2,9).Value  = "yes"
$ws.Cells.Item(2,10).Value = "Auto Testing 1"
$ws.Cells.Item(2,11).NumberFormat = "@"
$ws.Cells.Item(2,11).Value = "02/05/2017"
$ws.Cells.Item(2,12).Value = "10 minutes"

# 5. Row 3
$ws.Cells.Item(3,1).Value  = "Feedback"
$ws.Cells.Item(3,2).Value  = "Auto testing 4"
$ws.Cells.Item(3,3).Value  = "Medical Director"
$ws.Cells.Item(3,4).Value  = "National"
$ws.Cells.Item(3,5).Value  = "AutoTest 5 fo"
$ws.Cells.Item(3,6).Value  = "Auto Testing 1"
$ws.Cells.Item(3,7).Value  = "Academic"
$ws.Cells.Item(3,8).Value  = "yes"
$ws.Cells.Item(3,9).Value  = "yes"
$ws.Cells.Item(3,10).Value = "Auto Testing 2"
$ws.Cells.Item(3,11).NumberFormat = "@"
$ws.Cells.Item(3,11).Value = "02/05/2017"
$ws.Cells.Item(3,12).Value = "10 minutes"

# 6. Selection moves to H7
$ws.Range("H7").Select()
